# Update NATMI LR-pair TPM-derived statistics (Ncam1-Fgfr1) with the
# refreshed TPM numbers. Only the cells whose values actually move are
# touched; everything else (labels, styles, unrelated columns) is left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.04741066666666666
$ws.Range("H2").Value = 0.142232
$ws.Range("I2").Value = 0.003188134523263584
$ws.Range("J2").Value = 0.003188134523263585
$ws.Range("M2").Value = 2.294987
$ws.Range("N2").Value = 6.884961000000001
$ws.Range("O2").Value = 0.0158275801650097
$ws.Range("P2").Value = 0.0158275801650097
$ws.Range("Q2").Value = 0.1088068636613333
$ws.Range("R2").Value = 0.979261772952
$ws.Range("S2").Value = [double]"5.046045474378937E-05"
$ws.Range("T2").Value = [double]"5.046045474378937E-05"

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.04741066666666666
$ws.Range("H3").Value = 0.142232
$ws.Range("I3").Value = 0.003188134523263584
$ws.Range("J3").Value = 0.003188134523263585
$ws.Range("O3").Value = 0.769602070219672
$ws.Range("P3").Value = 0.7696020702196722
$ws.Range("Q3").Value = 5.290637397180443
$ws.Range("R3").Value = 47.61573657462399
$ws.Range("S3").Value = 0.002453594929242462
$ws.Range("T3").Value = 0.002453594929242462

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.04741066666666666
$ws.Range("H4").Value = 0.142232
$ws.Range("I4").Value = 0.003188134523263584
$ws.Range("J4").Value = 0.003188134523263585
$ws.Range("M4").Value = 31.11253633333333
$ws.Range("N4").Value = 93.337609
$ws.Range("O4").Value = 0.2145703496153182
$ws.Range("P4").Value = 0.2145703496153182
$ws.Range("Q4").Value = 1.475066089254222
$ws.Range("R4").Value = 13.275594803288
$ws.Range("S4").Value = 0.0006840791392773331
$ws.Range("T4").Value = 0.0006840791392773334

$ws.Range("I5").Value = 0.01595759596384214
$ws.Range("J5").Value = 0.01595759596384214
$ws.Range("M5").Value = 2.294987
$ws.Range("N5").Value = 6.884961000000001
$ws.Range("O5").Value = 0.0158275801650097
$ws.Range("P5").Value = 0.0158275801650097
$ws.Range("Q5").Value = 0.5446118900350001
$ws.Range("R5").Value = 4.901507010315
$ws.Range("S5").Value = 0.0002525701293585466
$ws.Range("T5").Value = 0.0002525701293585466

$ws.Range("I6").Value = 0.01595759596384214
$ws.Range("J6").Value = 0.01595759596384214
$ws.Range("O6").Value = 0.769602070219672
$ws.Range("P6").Value = 0.7696020702196722
$ws.Range("S6").Value = 0.01228099888950199
$ws.Range("T6").Value = 0.01228099888950199

$ws.Range("I7").Value = 0.01595759596384214
$ws.Range("J7").Value = 0.01595759596384214
$ws.Range("M7").Value = 31.11253633333333
$ws.Range("N7").Value = 93.337609
$ws.Range("O7").Value = 0.2145703496153182
$ws.Range("P7").Value = 0.2145703496153182
$ws.Range("Q7").Value = 7.383160434581667
$ws.Range("S7").Value = 0.003424026944981598
$ws.Range("T7").Value = 0.003424026944981598

$ws.Range("I8").Value = 0.9808542695128942
$ws.Range("J8").Value = 0.9808542695128943
$ws.Range("M8").Value = 2.294987
$ws.Range("N8").Value = 6.884961000000001
$ws.Range("O8").Value = 0.0158275801650097
$ws.Range("P8").Value = 0.0158275801650097
$ws.Range("Q8").Value = 33.47527401863734
$ws.Range("R8").Value = 301.277466167736
$ws.Range("S8").Value = 0.01552454958090736
$ws.Range("T8").Value = 0.01552454958090737

$ws.Range("I9").Value = 0.9808542695128942
$ws.Range("J9").Value = 0.9808542695128943
$ws.Range("O9").Value = 0.769602070219672
$ws.Range("P9").Value = 0.7696020702196722
$ws.Range("S9").Value = 0.7548674764009274
$ws.Range("T9").Value = 0.7548674764009278

$ws.Range("I10").Value = 0.9808542695128942
$ws.Range("J10").Value = 0.9808542695128943
$ws.Range("M10").Value = 31.11253633333333
$ws.Range("N10").Value = 93.337609
$ws.Range("O10").Value = 0.2145703496153182
$ws.Range("P10").Value = 0.2145703496153182
$ws.Range("S10").Value = 0.2104622435310593
$ws.Range("T10").Value = 0.2104622435310593

